$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.592.54'
$ws.Range("E2").Value = '  +6.60%  '

# Row 3
$ws.Range("D3").Value = '1.739.93'
$ws.Range("E3").Value = '  +4.60%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.18'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.88%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9971'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3741'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.74%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.56'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.95%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3402'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.204'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07521'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.97%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9978'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.458'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.09%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.48'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.77%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.073'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.46%  '

# Row 16
$ws.Range("D16").Value = '1.737.15'
$ws.Range("E16").Value = '  +4.30%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06702'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.53%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.05'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9955'
$ws.Range("D20").ClearFormats()

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.69%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.198'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.15'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +4.03%  '

# Row 24
$ws.Range("D24").Value = '26.552.83'
$ws.Range("E24").Value = '  +6.69%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.469'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.514'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.85%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.414'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +14.28%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.16'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.70'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.25%  '

# Row 30
$ws.Range("D30").Value = '1.930.39'
$ws.Range("E30").Value = '  +4.37%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '132.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.70%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.135'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.38%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.203'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08554'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.40%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.723'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.89%  '

# Row 36
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '13.17'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.85%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.486'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.57%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06361'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02355'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.66%  '

# Row 40
$ws.Range("E40").Value = '  +4.85%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.687'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.242'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6265'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.78%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.58'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +13.64%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9957'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.25%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.898'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.85%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6097'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +7.48%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.94'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.02%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.071'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +5.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07329'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.17%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.97'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.75%  '
